$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "A11" = -21.80360000000001
    "A12" = -21.5825
    "A15" = -21.89239999999999
    "A27" = -21.7256
    "A28" = -21.83990000000001
    "A31" = -21.7145
    "A32" = -21.31819999999998
    "A36" = -20.2417
    "A38" = -19.61799999999999
    "A46" = -21.90200000000001
    "A54" = -21.7043
    "A55" = -22.35520000000001
    "A56" = -22.1094
    "A67" = -21.48359999999999
    "A69" = -21.71199999999997
    "A72" = -21.7758
    "A73" = -19.8607
    "A83" = -21.7394
    "A86" = -21.9361
    "A91" = -21.50190000000001
    "A93" = -21.24309999999999
    "A99" = -19.99629999999999
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
